$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. A9 (GW_8 row): mark it as "completed" by copying the format
#        (green fill, centered) from A8, matching cellXfs s="11".
$ws.Range("A8").Copy() | Out-Null
$ws.Range("A9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 2. New row 10: GW_9 test case -----------------------------------
# Row height
$ws.Rows.Item(10).RowHeight = 112.5

# A10 -> "GW_9" (keeps the existing default style s="8")
$ws.Range("A10").Value = "GW_9"

# C10 -> "compatibility", reuse format+value from C9 (style s="7")
$ws.Range("C9").Copy($ws.Range("C10"))

# D10 -> "Main Search", reuse format+value from D4 (style s="9")
$ws.Range("D4").Copy($ws.Range("D10"))

# E10 -> rich-text description, new style: bold Calibri 11, left/top, wrap
$eCell = $ws.Range("E10")
$eCell.HorizontalAlignment = -4131
$eCell.VerticalAlignment = -4160
$eCell.WrapText = $true
$introText = "Check that the panel 'content info' displayed correctly  after opening the search results page:" + [char]10
$stepsText = "1. Open 'https://www.google.com/' with Chrome" + [char]10 + "2. Enter a valid query in the search (for example: wikipedia). " + [char]10 + "3. Use key 'Enter'. "
$eCell.Value = $introText + $stepsText
$eCell.Font.Name = "Calibri"
$eCell.Font.Size = 11
$eCell.Font.Bold = $true
$eCell.Font.ThemeColor = 1
$stepsChars = $eCell.Characters($introText.Length + 1, $stepsText.Length)
$stepsChars.Font.Name = "Calibri"
$stepsChars.Font.Size = 11
$stepsChars.Font.Bold = $false
$stepsChars.Font.ThemeColor = 1

# F10 -> Expected results text (keeps existing style s="5")
$ws.Range("F10").Value = "1. Website correctly open on Chrome browser." + [char]10 + "2. The entered text is displayed correctly in the input field." + [char]10 + "3. Google search page with query results has special panel 'content info' - footer for page"

# --- 3. Selection / scroll position, matching the target sheetView ---
$ws.Range("F13").Select() | Out-Null
